$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 345.6
$ws.Cells.Item(33, 9).Value = 307.55
$ws.Cells.Item(33, 10).Value = 421.7
$ws.Cells.Item(33, 11).Value = 307.55
$ws.Cells.Item(33, 12).Value = 421.7
$ws.Cells.Item(33, 13).Value = -78.55000000000001
$ws.Cells.Item(33, 14).Value = -879.7
$ws.Cells.Item(96, 8).Value = 991.7273
$ws.Cells.Item(96, 9).Value = 801.375
$ws.Cells.Item(96, 10).Value = 1499.3334
$ws.Cells.Item(96, 11).Value = 2404.125
$ws.Cells.Item(96, 12).Value = 4498.0002
$ws.Cells.Item(96, 13).Value = -1031.125
$ws.Cells.Item(96, 14).Value = -7244.0002
$ws.Cells.Item(103, 8).Value = 504.25
$ws.Cells.Item(103, 9).Value = 491
$ws.Cells.Item(103, 10).Value = 544
$ws.Cells.Item(103, 11).Value = 1473
$ws.Cells.Item(103, 12).Value = 1632
$ws.Cells.Item(103, 13).Value = -887
$ws.Cells.Item(103, 14).Value = -2804
$ws.Cells.Item(107, 8).Value = 3489.0715
$ws.Cells.Item(107, 9).Value = 1097.04
$ws.Cells.Item(107, 11).Value = 1097.04
$ws.Cells.Item(107, 13).Value = 822.96
$ws.Cells.Item(116, 8).Value = 5562643
$ws.Cells.Item(116, 9).Value = 8553156
$ws.Cells.Item(116, 10).Value = 8832.143
$ws.Cells.Item(116, 11).Value = 8553156
$ws.Cells.Item(116, 12).Value = 8832.143
$ws.Cells.Item(116, 13).Value = -8549714
$ws.Cells.Item(116, 14).Value = -15716.143
$ws.Cells.Item(138, 8).Value = 3659.982
$ws.Cells.Item(138, 9).Value = 1607.8572
$ws.Cells.Item(138, 10).Value = 4927.4707
$ws.Cells.Item(138, 11).Value = 4823.571599999999
$ws.Cells.Item(138, 12).Value = 14782.4121
$ws.Cells.Item(138, 13).Value = 316.4284000000007
$ws.Cells.Item(138, 14).Value = -25062.4121
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1822794
$ws.Cells.Item(32, 9).Value = 982.6111
$ws.Cells.Item(32, 11).Value = 982.6111
$ws.Cells.Item(32, 13).Value = -695.6111
$ws.Cells.Item(45, 8).Value = 1751.6842
$ws.Cells.Item(45, 9).Value = 1413.5
$ws.Cells.Item(45, 11).Value = 1413.5
$ws.Cells.Item(45, 13).Value = -1036.5
$ws.Cells.Item(61, 8).Value = 5109.391
$ws.Cells.Item(61, 9).Value = 8828.333000000001
$ws.Cells.Item(61, 11).Value = 8828.333000000001
$ws.Cells.Item(61, 13).Value = -8616.333000000001
$ws.Cells.Item(74, 8).Value = 2610.12
$ws.Cells.Item(74, 9).Value = 1936
$ws.Cells.Item(74, 11).Value = 1936
$ws.Cells.Item(74, 13).Value = -1062
$ws.Cells.Item(77, 8).Value = 2610.12
$ws.Cells.Item(77, 9).Value = 1936
$ws.Cells.Item(77, 11).Value = 9680
$ws.Cells.Item(77, 13).Value = -5312
$ws.Cells.Item(86, 8).Value = 103807.664
$ws.Cells.Item(86, 10).Value = 103807.664
$ws.Cells.Item(86, 12).Value = 103807.664
$ws.Cells.Item(86, 14).Value = -106179.664
$ws.Cells.Item(89, 8).Value = 103807.664
$ws.Cells.Item(89, 10).Value = 103807.664
$ws.Cells.Item(89, 12).Value = 311422.992
$ws.Cells.Item(89, 14).Value = -323278.992
$ws.Cells.Item(132, 8).Value = 643365
$ws.Cells.Item(132, 9).Value = 795898.4399999999
$ws.Cells.Item(132, 10).Value = 94244.60000000001
$ws.Cells.Item(132, 11).Value = 2387695.32
$ws.Cells.Item(132, 12).Value = 282733.8
$ws.Cells.Item(132, 13).Value = -2385165.32
$ws.Cells.Item(132, 14).Value = -287793.8
$ws.Cells.Item(136, 8).Value = 5109.391
$ws.Cells.Item(136, 9).Value = 8828.333000000001
$ws.Cells.Item(136, 11).Value = 26484.999
$ws.Cells.Item(136, 13).Value = -23934.999
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 20557
$ws.Cells.Item(75, 9).Value = 20557
$ws.Cells.Item(75, 11).Value = 20557
$ws.Cells.Item(75, 13).Value = -19621
$ws.Cells.Item(78, 8).Value = 20557
$ws.Cells.Item(78, 9).Value = 20557
$ws.Cells.Item(78, 11).Value = 61671
$ws.Cells.Item(78, 13).Value = -56991
$ws.Cells.Item(80, 8).Value = 11915050
$ws.Cells.Item(80, 10).Value = 17558128
$ws.Cells.Item(80, 12).Value = 17558128
$ws.Cells.Item(80, 14).Value = -17560124
$ws.Cells.Item(83, 8).Value = 11915050
$ws.Cells.Item(83, 10).Value = 17558128
$ws.Cells.Item(83, 12).Value = 87790640
$ws.Cells.Item(83, 14).Value = -87800624
$ws.Cells.Item(134, 8).Value = 1067973
$ws.Cells.Item(134, 9).Value = 1593346.6
$ws.Cells.Item(134, 10).Value = 17225.867
$ws.Cells.Item(134, 11).Value = 4780039.800000001
$ws.Cells.Item(134, 12).Value = 51677.601
$ws.Cells.Item(134, 13).Value = -4777504.800000001
$ws.Cells.Item(134, 14).Value = -56747.601
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2639.5881
$ws.Cells.Item(31, 9).Value = 998.5909
$ws.Cells.Item(31, 11).Value = 998.5909
$ws.Cells.Item(31, 13).Value = -703.5909
$ws.Cells.Item(34, 8).Value = 2639.5881
$ws.Cells.Item(34, 9).Value = 998.5909
$ws.Cells.Item(34, 11).Value = 998.5909
$ws.Cells.Item(34, 13).Value = -796.5909
$ws.Cells.Item(58, 8).Value = 41670870
$ws.Cells.Item(58, 9).Value = 66669916
$ws.Cells.Item(58, 10).Value = 5799.3335
$ws.Cells.Item(58, 11).Value = 66669916
$ws.Cells.Item(58, 12).Value = 5799.3335
$ws.Cells.Item(58, 13).Value = -66669713
$ws.Cells.Item(58, 14).Value = -6205.3335
$ws.Cells.Item(99, 8).Value = 10104870
$ws.Cells.Item(99, 10).Value = 3249
$ws.Cells.Item(99, 12).Value = 3249
$ws.Cells.Item(99, 14).Value = -6245
$ws.Cells.Item(105, 8).Value = 34485256
$ws.Cells.Item(105, 10).Value = 5492
$ws.Cells.Item(105, 12).Value = 5492
$ws.Cells.Item(105, 14).Value = -8986
$ws.Cells.Item(126, 8).Value = 10104870
$ws.Cells.Item(126, 10).Value = 3249
$ws.Cells.Item(126, 12).Value = 9747
$ws.Cells.Item(126, 14).Value = -14687
$ws.Cells.Item(136, 8).Value = 41670870
$ws.Cells.Item(136, 9).Value = 66669916
$ws.Cells.Item(136, 10).Value = 5799.3335
$ws.Cells.Item(136, 11).Value = 200009748
$ws.Cells.Item(136, 12).Value = 17398.0005
$ws.Cells.Item(136, 13).Value = -200007198
$ws.Cells.Item(136, 14).Value = -22498.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 3899.8647
$ws.Cells.Item(107, 10).Value = 4822.017
$ws.Cells.Item(107, 12).Value = 14466.051
$ws.Cells.Item(107, 14).Value = -18306.051
$ws.Cells.Item(121, 8).Value = 164327.25
$ws.Cells.Item(121, 10).Value = 218804.83
$ws.Cells.Item(121, 12).Value = 656414.49
$ws.Cells.Item(121, 14).Value = -659034.49
$ws.Cells.Item(137, 8).Value = 10891.9
$ws.Cells.Item(137, 9).Value = 8846.143
$ws.Cells.Item(137, 10).Value = 15665.333
$ws.Cells.Item(137, 11).Value = 26538.429
$ws.Cells.Item(137, 12).Value = 46995.999
$ws.Cells.Item(137, 13).Value = -21438.429
$ws.Cells.Item(137, 14).Value = -57195.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 55558548
$ws.Cells.Item(132, 9).Value = 66669360
$ws.Cells.Item(132, 11).Value = 200008080
$ws.Cells.Item(132, 13).Value = -200005550
$ws.Cells.Item(133, 8).Value = 96000
$ws.Cells.Item(133, 10).Value = 96000
$ws.Cells.Item(133, 12).Value = 96000
$ws.Cells.Item(133, 14).Value = -106120
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3798.7058
$ws.Cells.Item(40, 9).Value = 2314.923
$ws.Cells.Item(40, 11).Value = 2314.923
$ws.Cells.Item(40, 13).Value = -2178.923
$ws.Cells.Item(68, 8).Value = 2954.25
$ws.Cells.Item(68, 10).Value = 5000
$ws.Cells.Item(68, 12).Value = 5000
$ws.Cells.Item(68, 14).Value = -6498
$ws.Cells.Item(71, 8).Value = 2954.25
$ws.Cells.Item(71, 10).Value = 5000
$ws.Cells.Item(71, 12).Value = 25000
$ws.Cells.Item(71, 14).Value = -32488
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 14711247
$ws.Cells.Item(136, 9).Value = 17861360
$ws.Cells.Item(136, 10).Value = 10716.5
$ws.Cells.Item(136, 11).Value = 53584080
$ws.Cells.Item(136, 12).Value = 32149.5
$ws.Cells.Item(136, 13).Value = -53581530
$ws.Cells.Item(136, 14).Value = -37249.5
